{"js": "// Add the \"Customer\" list-item test cases under the existing bookmarked\n// (empty) numbered list paragraph. Mirrors the commit \"Added Add Customer\n// Tests\": the first bullet (which carries the _GoBack bookmark) gets the\n// text \"Customer null\", and 14 additional bullet paragraphs are appended\n// after it, reusing the same list formatting (NoSpacing style + numId 2).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text,isListItem\");\nawait context.sync();\n\n// Locate the (currently empty) numbered-list paragraph that holds the\n// _GoBack bookmark \u2014 it is the only list item in the document before the\n// edit is applied.\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.isListItem && para.text === \"\") {\n    anchor = para;\n    break;\n  }\n}\nif (!anchor) {\n  throw new Error(\"Could not locate the target list paragraph (_GoBack bookmark item).\");\n}\n\n// New bullet texts to append, in order. Entries with a second element\n// get that text appended as extra content on the same paragraph.\nconst newItems = [\n  [\"Customer cu FirstName de 1 caracter\"],\n  [\"Customer cu FirstName de 31 de caractere\"],\n  [\"Customer doar cu campurile required, pentru a invalida regex-ul de pe campurile care nu sunt mandatory\"],\n  [\"Customer fara campurile mandatorii (fara CNP)\"],\n  [\"Customer cu Email corect\"],\n  [\"Customer cu Email gresit\"],\n  [\"Customer cu Email gresit 2\"],\n  [\"Customer cu telefon corect\", \" fara caractere speciale\"],\n  [\"Customer cu telefon corect si caractere speciale\", \" (-)\"],\n  [\"Customer cu telefon corect si caractere speciale (.)\"],\n  [\"Customer cu telefon incorect (prea scurt)\"],\n  [\"Customer cu telefon incorect\", \" (nu contine doar cifre)\"],\n  [\"Customer cu telefon \"],\n  [\"Customer complet\"],\n];\n\n// Fill the bookmarked paragraph with its text first.\nanchor.insertText(\"Customer null\", Word.InsertLocation.end);\n\n// Insert the remaining bullets right after the anchor, each one becoming\n// the new \"previous\" paragraph so the list stays in order.\nlet previous = anchor;\nfor (const item of newItems) {\n  const newPara = previous.insertParagraph(item[0], Word.InsertLocation.after);\n  if (item.length > 1) {\n    newPara.insertText(item[1], Word.InsertLocation.end);\n  }\n  previous = newPara;\n}\n\nawait context.sync();\n", "ps1": "# Add the \"Customer\" list-item test cases under the existing bookmarked\n# (empty) numbered list paragraph. Mirrors the commit \"Added Add Customer\n# Tests\": the first bullet (which carries the _GoBack bookmark) gets the\n# text \"Customer null\", and 14 additional bullet paragraphs are appended\n# after it, reusing the same list formatting (NoSpacing style + numId 2).\n\n$doc = $word.ActiveDocument\n\n# Locate the index of the (currently empty) numbered-list paragraph that\n# holds the _GoBack bookmark -- it is the only list item in the document\n# before the edit is applied.\n$anchorIndex = 0\n$idx = 0\nforeach ($p in $doc.Paragraphs) {\n    $idx = $idx + 1\n    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)\n    $listType = $p.Range.ListFormat.ListType\n    if ($listType -ne 0 -and $txt -eq \"\") {\n        $anchorIndex = $idx\n    }\n}\nif ($anchorIndex -eq 0) {\n    throw \"Could not locate the target list paragraph (_GoBack bookmark item).\"\n}\n\n# Fill the bookmarked paragraph with its text first.\n$anchor = $doc.Paragraphs($anchorIndex)\n$anchor.Range.InsertAfter(\"Customer null\")\n\n# New bullet texts to append, in order. Entries with a second element\n# get that text appended as extra content on the same paragraph.\n$items = @(\n    @(\"Customer cu FirstName de 1 caracter\"),\n    @(\"Customer cu FirstName de 31 de caractere\"),\n    @(\"Customer doar cu campurile required, pentru a invalida regex-ul de pe campurile care nu sunt mandatory\"),\n    @(\"Customer fara campurile mandatorii (fara CNP)\"),\n    @(\"Customer cu Email corect\"),\n    @(\"Customer cu Email gresit\"),\n    @(\"Customer cu Email gresit 2\"),\n    @(\"Customer cu telefon corect\", \" fara caractere speciale\"),\n    @(\"Customer cu telefon corect si caractere speciale\", \" (-)\"),\n    @(\"Customer cu telefon corect si caractere speciale (.)\"),\n    @(\"Customer cu telefon incorect (prea scurt)\"),\n    @(\"Customer cu telefon incorect\", \" (nu contine doar cifre)\"),\n    @(\"Customer cu telefon \"),\n    @(\"Customer complet\")\n)\n\n# Insert the remaining bullets right after the anchor paragraph, each one\n# becoming the new \"current\" paragraph so the list stays in order.\n$currentIndex = $anchorIndex\nforeach ($item in $items) {\n    $doc.Paragraphs($currentIndex).Range.InsertParagraphAfter()\n    $currentIndex = $currentIndex + 1\n    $newPara = $doc.Paragraphs($currentIndex)\n    $newPara.Range.InsertAfter($item[0])\n    if ($item.Count -gt 1) {\n        $doc.Paragraphs($currentIndex).Range.InsertAfter($item[1])\n    }\n}\n"}
